# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with freshly scraped values (GitHub Actions refresh).
#
# Note: Price values such as "312.08" or "1.005" are valid-looking numbers,
# but the sheet stores them as text (thousand separators use '.', so values
# like "27.065.72" are not numeric). To stop Excel from auto-converting
# assignments like "312.08" into a numeric cell, we prefix the literal with
# a leading apostrophe (forces text entry) and then reset the cell's style
# back to "Normal" so no stray quote-prefix formatting lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'27.065.72"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "
$cell = $ws.Range("D3")
$cell.Value = "'1.823.86"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.30%  "
$cell = $ws.Range("D5")
$cell.Value = "'312.08"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$cell = $ws.Range("D6")
$cell.Value = "'1.005"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "
$cell = $ws.Range("D7")
$cell.Value = "'0.4690"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$cell = $ws.Range("D8")
$cell.Value = "'0.3652"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.27%  "
$cell = $ws.Range("D9")
$cell.Value = "'0.07379"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.25%  "
$cell = $ws.Range("D10")
$cell.Value = "'0.8774"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.35%  "
$cell = $ws.Range("D11")
$cell.Value = "'20.23"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.40%  "
$cell = $ws.Range("D12")
$cell.Value = "'1.888.44"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +2.20%  "
$cell = $ws.Range("D13")
$cell.Value = "'0.07597"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +4.56%  "
$cell = $ws.Range("D14")
$cell.Value = "'5.364"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.24%  "
$cell = $ws.Range("D15")
$cell.Value = "'92.95"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "
$cell = $ws.Range("D16")
$cell.Value = "'6.522"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "
$cell = $ws.Range("D17")
$cell.Value = "'1.004"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "
$cell = $ws.Range("D18")
$cell.Value = "'0.000008710"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$cell = $ws.Range("D20")
$cell.Value = "'27.432.43"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.54%  "
$cell = $ws.Range("D21")
$cell.Value = "'14.58"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.76%  "
$cell = $ws.Range("D22")
$cell.Value = "'5.232"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.10%  "
$cell = $ws.Range("D23")
$cell.Value = "'10.61"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$cell = $ws.Range("D24")
$cell.Value = "'2.082.05"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.73%  "
$cell = $ws.Range("D25")
$cell.Value = "'1.878"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "
$cell = $ws.Range("D26")
$cell.Value = "'151.32"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.29%  "
$cell = $ws.Range("D27")
$cell.Value = "'18.35"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "
$cell = $ws.Range("D28")
$cell.Value = "'2.128"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.84%  "
$cell = $ws.Range("D29")
$cell.Value = "'5.168"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.57%  "
$cell = $ws.Range("D30")
$cell.Value = "'116.29"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.59%  "
$cell = $ws.Range("D31")
$cell.Value = "'0.08911"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.28%  "
$cell = $ws.Range("D32")
$cell.Value = "'0.7437"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.45%  "
$cell = $ws.Range("D33")
$cell.Value = "'1.159"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.36%  "
$cell = $ws.Range("D34")
$cell.Value = "'4.508"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +0.29%  "
$cell = $ws.Range("D36")
$cell.Value = "'2.730"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +14.65%  "
$ws.Range("E37").Value = "  +0.27%  "
$cell = $ws.Range("D38")
$cell.Value = "'1.086"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.85%  "
$cell = $ws.Range("D39")
$cell.Value = "'0.05288"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.47%  "
$cell = $ws.Range("D40")
$cell.Value = "'0.01930"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "
$cell = $ws.Range("D41")
$cell.Value = "'7.295"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("E42").Value = "  -1.62%  "
$cell = $ws.Range("D43")
$cell.Value = "'0.5256"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.96%  "
$cell = $ws.Range("D44")
$cell.Value = "'0.1641"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.83%  "
$cell = $ws.Range("D45")
$cell.Value = "'8.358"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.46%  "
$cell = $ws.Range("D46")
$cell.Value = "'0.4899"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  -1.64%  "
$cell = $ws.Range("D48")
$cell.Value = "'1.005"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  +0.99%  "
$cell = $ws.Range("D50")
$cell.Value = "'1.650"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "
$cell = $ws.Range("D51")
$cell.Value = "'0.06262"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.61%  "
